$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 17: "proceso migracion datos" task now has a responsible person
# (Agustina) and its status column filled in ("en proceso"); the status
# cell was previously just an empty, underline-formatted placeholder, so
# drop that leftover underline once it holds real text.
$ws.Range("B17").Value = "Agustina"
$ws.Range("C17").Value = "en proceso"
$ws.Range("C17").Font.Underline = $false

# Row 31: fill in the status column with "en proceso"
$ws.Range("C31").Value = "en proceso"

# Row 40 ("acentos!"): assign Lucas as responsible and mark 100% complete
$ws.Range("B40").Value = "Lucas"
$ws.Range("C40").Value = 1
$ws.Range("C40").NumberFormat = "0%"

# Reset the view: scroll position / selection back to A19 (no more
# multi-cell selection further down the sheet)
$ws.Activate()
$ws.Range("A19").Select()
